$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing amounts (column B, rows 2-5)
$ws.Range("B2").Value = 4534504.109000016
$ws.Range("B3").Value = 1958333.170000002
$ws.Range("B4").Value = 274656.1250000001
$ws.Range("B5").Value = 1811653

# Add new row 6: municipal / 1107319 (copy formatting from the row above)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "municipal"
$ws.Range("B6").Value = 1107319
